# Normalize the "Recorded By" (column G) attendance-recorder lists so the
# "System" / admin entries sort consistently, e.g.
#   "dnasr281@gmail.com, System"                  -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System"                  -> "System, backup@backdoor.com"
#   "system, backup@backdoor.com, System"          -> "system, System, backup@backdoor.com"
#   "dnasr281@gmail.com, admin@admin.com"          -> "admin@admin.com, dnasr281@gmail.com"
#
# Only cells whose text matches one of these exact combinations are touched;
# everything else in the "Recorded By" column (single recorder, or other
# combinations) is left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$lastRow = $ws.UsedRange.Rows.Count
$col = "G"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("$col$r")
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
